$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header text updates (report volume/number and week-covering dates) ----
$ws.Range("A8").Value = "Volume 32   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  6/30/2025  Through  7/6/2025"

# ---- Crime-statistics table updates ----
# Helper reference cells whose styles we reuse when a cell's underlying
# type needs to flip between "text placeholder" and "number":
#   C22 -> general/text style (s=13)
#   I29 -> integer/count style (s=14)
#   L29 -> percentage style (s=15)

# Row 15 (Rape)
$ws.Range("D15").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("N15").Value = -57.142857142857

# Row 16 (Robbery)
$ws.Range("C16").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -100
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -20
$ws.Range("J16").Value = 31
$ws.Range("K16").Value = -51.612903225806
$ws.Range("L16").Value = -57.142857142857
$ws.Range("N16").Value = -91.017964071856

# Row 17 (Fel. Assault)
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 42.857142857142
$ws.Range("I17").Value = 65
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 25
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 140.740740740741
$ws.Range("N17").Value = 0

# Row 18 (Burglary)
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -16.666666666666
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = 29.411764705882
$ws.Range("I18").Value = 151
$ws.Range("J18").Value = 133
$ws.Range("K18").Value = 13.533834586466
$ws.Range("L18").Value = -8.484848484848
$ws.Range("M18").Value = 22.764227642276
$ws.Range("N18").Value = -71.509433962264

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 54.838709677419
$ws.Range("I19").Value = 228
$ws.Range("J19").Value = 239
$ws.Range("K19").Value = -4.602510460251
$ws.Range("L19").Value = -33.137829912023
$ws.Range("M19").Value = 28.089887640449
$ws.Range("N19").Value = -13.962264150943

# Row 20 (G.L.A.)
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = 85.714285714285
$ws.Range("I20").Value = 134
$ws.Range("J20").Value = 128
$ws.Range("K20").Value = 4.6875
$ws.Range("L20").Value = 78.666666666666
$ws.Range("M20").Value = 83.561643835616
$ws.Range("N20").Value = -91.956782713085

# Row 21 (TOTAL)
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 5.555555555555
$ws.Range("F21").Value = 110
$ws.Range("G21").Value = 75
$ws.Range("H21").Value = 46.666666666666
$ws.Range("I21").Value = 598
$ws.Range("J21").Value = 587
$ws.Range("K21").Value = 1.873935264054
$ws.Range("L21").Value = -11.538461538461
$ws.Range("M21").Value = 33.184855233853
$ws.Range("N21").Value = -77.860051832654

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 51
$ws.Range("G24").Value = 41
$ws.Range("H24").Value = 24.390243902439
$ws.Range("I24").Value = 280
$ws.Range("J24").Value = 286
$ws.Range("K24").Value = -2.097902097902
$ws.Range("L24").Value = -4.436860068259
$ws.Range("M24").Value = 14.754098360655

# Row 25 (Retail Theft)
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = -36.363636363636
$ws.Range("I25").Value = 52
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = -13.333333333333
$ws.Range("L25").Value = -22.388059701492

# Row 26 (Misd. Assault)
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 12.5
$ws.Range("F26").Value = 20
$ws.Range("G26").Value = 23
$ws.Range("H26").Value = -13.043478260869
$ws.Range("I26").Value = 103
$ws.Range("J26").Value = 109
$ws.Range("K26").Value = -5.504587155963
$ws.Range("L26").Value = -8.035714285714
$ws.Range("M26").Value = 25.609756097561

# Row 27 (UCR Rape*)
$ws.Range("D27").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# Row 28 (Other Sex Crimes)
$ws.Range("C28").Value = 1
$ws.Range("I29").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 5
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = -37.5

# Row 31 (Hate Crimes)
$ws.Range("D31").Value = 3
$ws.Range("I29").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").Value = -100
$ws.Range("L29").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("G31").Value = 3
$ws.Range("I29").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("H31").Value = -100
$ws.Range("L29").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = -16.666666666666

# Row 33 (Traffic Fatalities)
$ws.Range("D33").Value = "'0"
$ws.Range("C22").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").Value = "'***.*"
$ws.Range("C22").Copy()
$ws.Range("E33").PasteSpecial(-4122)
